$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C143 previously held the literal text "NA"; the updated export leaves
# it an empty (inline-string-typed) cell instead.
$ws.Range("C143").Value = ""

# Two new rows appended by the refreshed script run. The "Date" column
# values look like dates, so Excel's autoconvert would otherwise turn
# them into date serials -- force them to stay plain text, matching the
# rest of the sheet, then drop back to the default style so no extra
# number-format style sticks to the cell.
$ws.Range("A144").NumberFormat = "@"
$ws.Range("A144").Value = "2025-06-27"
$ws.Range("A144").Style = "Normal"
$ws.Range("B144").Value = "eaux de surface"
$ws.Range("C144").Value = 113
$ws.Range("D144").Value = 1

$ws.Range("A145").NumberFormat = "@"
$ws.Range("A145").Value = "2025-06-27"
$ws.Range("A145").Style = "Normal"
$ws.Range("B145").Value = "eaux de surface"
$ws.Range("C145").Value = 118
$ws.Range("D145").Value = 1
